$wb = $excel.ActiveWorkbook

# --- RocketMarket sheet ---
$rocket = $wb.Worksheets.Item("RocketMarket")

# Row 7 (Salvo Rocket): destroy_up_to_2_ships -> destroy_up_to_2_ships_then_lose_one_1_bank_currency
$rocket.Range("E7").Value = "destroy_up_to_2_ships_then_lose_one_1_bank_currency"
$rocket.Range("F7").Value = "Two hits, then lose 1 bank currency"

# Row 10 (Twin Salvo): destroy_up_to_2_ships -> destory_up_to_2_ships
$rocket.Range("E10").Value = "destory_up_to_2_ships"

# --- ShieldMarket sheet ---
$shield = $wb.Worksheets.Item("ShieldMarket")

# Row 3 (Decoy Drone): assign_to_ship_block_1 -> assign_to_ship_block_1_draw_1_discard_1
$shield.Range("E3").Value = "assign_to_ship_block_1_draw_1_discard_1"
$shield.Range("F3").Value = "Assign: block 1; immediately draw 1 then discard 1"
